$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices for C, E..Y (skip D) that changed in this edit
$cols = @(3, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25)

$data = @{
    2 = @(1.131002426147461, 2285.848009800326, 0.1492674433751622, 0.09927312128244237, 0.08356394355512026, 0.07094213120565727, 0.06374305997962727, 0.06042746198794464, 0.05630360086731506, 0.05376523078475696, 0.05180309660096867, 0.04960818795206661, 0.04901414518779094, 0.04837361249679491, 0.04731814764838346, 0.04630324371111674, 0.04591971167623956, 0.04534233822172413, 0.0451183433964567, 0.04463444409226793, 0.04463444409226793, 0.04455844073684844)
    3 = @(1.190001726150513, 2362.986396963865, 0.142518954151997, 0.104924393502459, 0.08605364978551151, 0.07472330671174886, 0.06799877980648114, 0.06407211526365073, 0.06046114218799448, 0.0570009633516655, 0.05399215083133196, 0.05141699583495936, 0.05045752684389462, 0.04865556793157465, 0.04783956815136844, 0.04754354116390101, 0.04711166554577965, 0.04677639768616505, 0.04654384518588429, 0.04628456918939834, 0.04613053831405734, 0.04606211300124492)
    4 = @(1.260991811752319, 2383.158573206422, 0.1597815532138741, 0.1253258001093468, 0.1014737287122646, 0.0841952438886662, 0.07286675371608217, 0.06837339469812451, 0.06133123828773693, 0.05722132850364586, 0.05554737553449812, 0.05386505168302744, 0.05134690024041998, 0.05061573506943681, 0.04877636775381985, 0.04819422613550696, 0.04762666636509488, 0.04746910083031759, 0.04720758485143456, 0.04679346753995573, 0.04654805430447385, 0.04645533281104136)
    5 = @(1.094999551773071, 2468.072067458463, 0.1560211332194928, 0.1139715125125783, 0.09425086627053438, 0.08227598158509403, 0.07341310455469066, 0.06866760843607088, 0.0642129113050203, 0.06124358251091944, 0.05805259529311178, 0.05573763543145198, 0.05371759157563079, 0.05240420446385861, 0.05053398558361499, 0.04988151982314276, 0.04933834129776322, 0.0491608504091707, 0.04862180270722717, 0.04843348271426458, 0.04811056661712402, 0.04811056661712402)
    6 = @(1.08699893951416, 2265.73661717935, 0.1569513906080221, 0.1145969035750158, 0.0878055940856502, 0.07751555953865025, 0.06392747071487516, 0.06280885706358273, 0.05690542300489374, 0.05426424343776942, 0.05085165263168894, 0.0491783133956552, 0.04821242079287541, 0.0474903587538722, 0.04654083827867433, 0.04535546648210063, 0.04527476729455628, 0.04492010272346388, 0.04463757835488501, 0.04439264935692733, 0.04420694926437773, 0.04416640579296978)
    7 = @(1.157040119171143, 2352.684319576196, 0.1550122527548267, 0.1100146301376521, 0.08382014488319389, 0.07331600036465011, 0.06808463662108849, 0.06189364923249952, 0.05909631201699142, 0.05623909263071907, 0.05369488767192627, 0.05196913272154202, 0.05019685912516, 0.04846097577208427, 0.04784620632332751, 0.04744387233475768, 0.04714401845024079, 0.04647603446757589, 0.0464513904704754, 0.04597399509160776, 0.04597399509160776, 0.04586129277926307)
    8 = @(1.116955518722534, 2411.399905805347, 0.1581735178394743, 0.1151072765022438, 0.08758054339072222, 0.07763434698001426, 0.06915994076357546, 0.06263384221932879, 0.05834116406380992, 0.05589758898682473, 0.05399699699508318, 0.05203118666239887, 0.05072641937310073, 0.04940227929220654, 0.04879686319712696, 0.04853150123953228, 0.04799773772654729, 0.04779700662334977, 0.04736408415554227, 0.04716691497779268, 0.04702763746077111, 0.04700584611706328)
    9 = @(1.317004203796387, 2276.414977535259, 0.1576915360629395, 0.1128895512954495, 0.0934351732615671, 0.07809099820825018, 0.0697054198824677, 0.06086284299627368, 0.0575400081423296, 0.05323320524962289, 0.04903408897605753, 0.04903408897605753, 0.04838660241012335, 0.04754397832309903, 0.04698481320790254, 0.04608641376181987, 0.04555903653330343, 0.04522077187801338, 0.044979784835726, 0.04452432437335015, 0.04452432437335015, 0.04437456096559958)
    10 = @(1.1629958152771, 2290.704095957106, 0.155347510747666, 0.1164104181492145, 0.08649217211694712, 0.07296873653094224, 0.06384972792215189, 0.05905698642189849, 0.05394857392454826, 0.052710317435791, 0.05012317142428308, 0.04744869778508713, 0.04720219055050031, 0.04711515187910845, 0.046011653304181, 0.046011653304181, 0.04582807128646074, 0.04539591579879913, 0.04501762408629759, 0.04474130201131095, 0.04473971006833578, 0.04465310128571356)
    11 = @(1.191977739334106, 2300.205917159951, 0.1577255335033474, 0.1171180033527798, 0.09452084297694177, 0.07647482893275717, 0.0610476541669706, 0.06008237746882798, 0.05468951491541064, 0.05292749173776914, 0.05143223255301915, 0.04966243436083274, 0.04850126991900149, 0.04776537280549033, 0.04686487193968189, 0.04636640177525734, 0.04572593831358305, 0.04533756330361025, 0.04515364188010709, 0.04501804778834525, 0.04501804778834525, 0.04483832197192886)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item([int]$row, $cols[$i]).Value = $vals[$i]
    }
}